$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 66.47695399999999
$ws.Range("H2").Value = 199.430862
$ws.Range("I2").Value = 0.04311983106164722
$ws.Range("J2").Value = 0.04311983106164721
$ws.Range("M2").Value = 6.066605666666667
$ws.Range("N2").Value = 18.199817
$ws.Range("O2").Value = 0.8497846287916651
$ws.Range("P2").Value = 0.8497846287916652
$ws.Range("Q2").Value = 403.2894658391393
$ws.Range("R2").Value = 3629.605192552254
$ws.Range("S2").Value = 0.03664256963228119
$ws.Range("T2").Value = 0.03664256963228119
$ws.Range("G3").Value = 66.47695399999999
$ws.Range("H3").Value = 199.430862
$ws.Range("I3").Value = 0.04311983106164722
$ws.Range("J3").Value = 0.04311983106164721
$ws.Range("O3").Value = 0.1196497582104962
$ws.Range("P3").Value = 0.1196497582104962
$ws.Range("Q3").Value = 56.78319593178131
$ws.Range("R3").Value = 511.0487633860319
$ws.Range("S3").Value = 0.005159277360603533
$ws.Range("T3").Value = 0.005159277360603533
$ws.Range("G4").Value = 66.47695399999999
$ws.Range("H4").Value = 199.430862
$ws.Range("I4").Value = 0.04311983106164722
$ws.Range("J4").Value = 0.04311983106164721
$ws.Range("M4").Value = 0.1824346666666667
$ws.Range("N4").Value = 0.547304
$ws.Range("O4").Value = 0.02555468148257719
$ws.Range("P4").Value = 0.02555468148257719
$ws.Range("Q4").Value = 12.12770094400533
$ws.Range("R4").Value = 109.149308496048
$ws.Range("S4").Value = 0.001101913548362933
$ws.Range("T4").Value = 0.001101913548362933
$ws.Range("G5").Value = 66.47695399999999
$ws.Range("H5").Value = 199.430862
$ws.Range("I5").Value = 0.04311983106164722
$ws.Range("J5").Value = 0.04311983106164721
$ws.Range("M5").Value = 0.035773
$ws.Range("N5").Value = 0.107319
$ws.Range("O5").Value = 0.005010931515261538
$ws.Range("P5").Value = 0.005010931515261539
$ws.Range("Q5").Value = 2.378080075442
$ws.Range("R5").Value = 21.402720678978
$ws.Range("S5").Value = 0.0002160705203995614
$ws.Range("T5").Value = 0.0002160705203995614
$ws.Range("I6").Value = 0.8830494168872806
$ws.Range("J6").Value = 0.8830494168872804
$ws.Range("M6").Value = 6.066605666666667
$ws.Range("N6").Value = 18.199817
$ws.Range("O6").Value = 0.8497846287916651
$ws.Range("P6").Value = 0.8497846287916652
$ws.Range("Q6").Value = 8258.949974476791
$ws.Range("R6").Value = 74330.54977029111
$ws.Range("S6").Value = 0.750401820934254
$ws.Range("T6").Value = 0.750401820934254
$ws.Range("I7").Value = 0.8830494168872806
$ws.Range("J7").Value = 0.8830494168872804
$ws.Range("O7").Value = 0.1196497582104962
$ws.Range("P7").Value = 0.1196497582104962
$ws.Range("S7").Value = 0.1056566492184828
$ws.Range("T7").Value = 0.1056566492184828
$ws.Range("I8").Value = 0.8830494168872806
$ws.Range("J8").Value = 0.8830494168872804
$ws.Range("M8").Value = 0.1824346666666667
$ws.Range("N8").Value = 0.547304
$ws.Range("O8").Value = 0.02555468148257719
$ws.Range("P8").Value = 0.02555468148257719
$ws.Range("Q8").Value = 248.362736659992
$ws.Range("R8").Value = 2235.264629939928
$ws.Range("S8").Value = 0.02256604658192997
$ws.Range("T8").Value = 0.02256604658192997
$ws.Range("I9").Value = 0.8830494168872806
$ws.Range("J9").Value = 0.8830494168872804
$ws.Range("M9").Value = 0.035773
$ws.Range("N9").Value = 0.107319
$ws.Range("O9").Value = 0.005010931515261538
$ws.Range("P9").Value = 0.005010931515261539
$ws.Range("Q9").Value = 48.70061343533699
$ws.Range("R9").Value = 438.305520918033
$ws.Range("S9").Value = 0.004424900152613798
$ws.Range("T9").Value = 0.004424900152613798
$ws.Range("G10").Value = 44.831112
$ws.Range("H10").Value = 134.493336
$ws.Range("I10").Value = 0.02907940059566787
$ws.Range("J10").Value = 0.02907940059566786
$ws.Range("M10").Value = 6.066605666666667
$ws.Range("N10").Value = 18.199817
$ws.Range("O10").Value = 0.8497846287916651
$ws.Range("P10").Value = 0.8497846287916652
$ws.Range("Q10").Value = 271.972678102168
$ws.Range("R10").Value = 2447.754102919512
$ws.Range("S10").Value = 0.02471122764067375
$ws.Range("T10").Value = 0.02471122764067374
$ws.Range("G11").Value = 44.831112
$ws.Range("H11").Value = 134.493336
$ws.Range("I11").Value = 0.02907940059566787
$ws.Range("J11").Value = 0.02907940059566786
$ws.Range("O11").Value = 0.1196497582104962
$ws.Range("P11").Value = 0.1196497582104962
$ws.Range("Q11").Value = 38.29377947334399
$ws.Range("R11").Value = 344.644015260096
$ws.Range("S11").Value = 0.003479343250177819
$ws.Range("T11").Value = 0.003479343250177819
$ws.Range("G12").Value = 44.831112
$ws.Range("H12").Value = 134.493336
$ws.Range("I12").Value = 0.02907940059566787
$ws.Range("J12").Value = 0.02907940059566786
$ws.Range("M12").Value = 0.1824346666666667
$ws.Range("N12").Value = 0.547304
$ws.Range("O12").Value = 0.02555468148257719
$ws.Range("P12").Value = 0.02555468148257719
$ws.Range("Q12").Value = 8.178748974015999
$ws.Range("R12").Value = 73.608740766144
$ws.Range("S12").Value = 0.0007431148199265578
$ws.Range("T12").Value = 0.0007431148199265577
$ws.Range("G13").Value = 44.831112
$ws.Range("H13").Value = 134.493336
$ws.Range("I13").Value = 0.02907940059566787
$ws.Range("J13").Value = 0.02907940059566786
$ws.Range("M13").Value = 0.035773
$ws.Range("N13").Value = 0.107319
$ws.Range("O13").Value = 0.005010931515261538
$ws.Range("P13").Value = 0.005010931515261539
$ws.Range("Q13").Value = 1.603743369576
$ws.Range("R13").Value = 14.433690326184
$ws.Range("S13").Value = 0.0001457148848897473
$ws.Range("T13").Value = 0.0001457148848897473
$ws.Range("G14").Value = 52.83062100000001
$ws.Range("H14").Value = 158.491863
$ws.Range("I14").Value = 0.0342682285413064
$ws.Range("J14").Value = 0.03426822854130639
$ws.Range("M14").Value = 6.066605666666667
$ws.Range("N14").Value = 18.199817
$ws.Range("O14").Value = 0.8497846287916651
$ws.Range("P14").Value = 0.8497846287916652
$ws.Range("Q14").Value = 320.5025447321191
$ws.Range("R14").Value = 2884.522902589071
$ws.Range("S14").Value = 0.029120613870322
$ws.Range("T14").Value = 0.029120613870322
$ws.Range("G15").Value = 52.83062100000001
$ws.Range("H15").Value = 158.491863
$ws.Range("I15").Value = 0.0342682285413064
$ws.Range("J15").Value = 0.03426822854130639
$ws.Range("O15").Value = 0.1196497582104962
$ws.Range("P15").Value = 0.1196497582104962
$ws.Range("Q15").Value = 45.126789404952
$ws.Range("R15").Value = 406.141104644568
$ws.Range("S15").Value = 0.004100185259269334
$ws.Range("T15").Value = 0.004100185259269334
$ws.Range("G16").Value = 52.83062100000001
$ws.Range("H16").Value = 158.491863
$ws.Range("I16").Value = 0.0342682285413064
$ws.Range("J16").Value = 0.03426822854130639
$ws.Range("M16").Value = 0.1824346666666667
$ws.Range("N16").Value = 0.547304
$ws.Range("O16").Value = 0.02555468148257719
$ws.Range("P16").Value = 0.02555468148257719
$ws.Range("Q16").Value = 9.638136731928002
$ws.Range("R16").Value = 86.74323058735202
$ws.Range("S16").Value = 0.0008757136653452457
$ws.Range("T16").Value = 0.0008757136653452457
$ws.Range("G17").Value = 52.83062100000001
$ws.Range("H17").Value = 158.491863
$ws.Range("I17").Value = 0.0342682285413064
$ws.Range("J17").Value = 0.03426822854130639
$ws.Range("M17").Value = 0.035773
$ws.Range("N17").Value = 0.107319
$ws.Range("O17").Value = 0.005010931515261538
$ws.Range("P17").Value = 0.005010931515261539
$ws.Range("Q17").Value = 1.889909805033
$ws.Range("R17").Value = 17.009188245297
$ws.Range("S17").Value = 0.0001717157463698172
$ws.Range("T17").Value = 0.0001717157463698172
$ws.Range("G18").Value = 16.16161433333333
$ws.Range("H18").Value = 48.484843
$ws.Range("I18").Value = 0.01048312291409786
$ws.Range("J18").Value = 0.01048312291409786
$ws.Range("M18").Value = 6.066605666666667
$ws.Range("N18").Value = 18.199817
$ws.Range("O18").Value = 0.8497846287916651
$ws.Range("P18").Value = 0.8497846287916652
$ws.Range("Q18").Value = 98.04614109708122
$ws.Range("R18").Value = 882.4152698737309
$ws.Range("S18").Value = 0.008908396714134051
$ws.Range("T18").Value = 0.008908396714134051
$ws.Range("G19").Value = 16.16161433333333
$ws.Range("H19").Value = 48.484843
$ws.Range("I19").Value = 0.01048312291409786
$ws.Range("J19").Value = 0.01048312291409786
$ws.Range("O19").Value = 0.1196497582104962
$ws.Range("P19").Value = 0.1196497582104962
$ws.Range("Q19").Value = 13.80490618242755
$ws.Range("R19").Value = 124.244155641848
$ws.Range("S19").Value = 0.001254303121962721
$ws.Range("T19").Value = 0.001254303121962722
$ws.Range("G20").Value = 16.16161433333333
$ws.Range("H20").Value = 48.484843
$ws.Range("I20").Value = 0.01048312291409786
$ws.Range("J20").Value = 0.01048312291409786
$ws.Range("M20").Value = 0.1824346666666667
$ws.Range("N20").Value = 0.547304
$ws.Range("O20").Value = 0.02555468148257719
$ws.Range("P20").Value = 0.02555468148257719
$ws.Range("Q20").Value = 2.948438723696889
$ws.Range("R20").Value = 26.535948513272
$ws.Range("S20").Value = 0.0002678928670124773
$ws.Range("T20").Value = 0.0002678928670124773
$ws.Range("G21").Value = 16.16161433333333
$ws.Range("H21").Value = 48.484843
$ws.Range("I21").Value = 0.01048312291409786
$ws.Range("J21").Value = 0.01048312291409786
$ws.Range("M21").Value = 0.035773
$ws.Range("N21").Value = 0.107319
$ws.Range("O21").Value = 0.005010931515261538
$ws.Range("P21").Value = 0.005010931515261539
$ws.Range("Q21").Value = 0.5781494295463333
$ws.Range("R21").Value = 5.203344865917
$ws.Range("S21").Value = 0.00005253021098861336
$ws.Range("T21").Value = 0.00005253021098861336
